$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.525.55"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.853.59"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.77"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4700"
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2745"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06343"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.49"
$ws.Range("E10").Value = "  +6.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.859.82"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07436"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.089"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.63"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6276"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.508.99"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "243.57"
$ws.Range("E17").Value = "  +5.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.71"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007344"
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.963"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.006"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.283"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.11"
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.05"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.884"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1013"
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.371"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.042"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.854"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04902"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7065"
$ws.Range("E34").Value = "  -2.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.706"
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01906"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.685"
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8759"
$ws.Range("E38").Value = "  -3.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.981"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "105.04"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4073"
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.510"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.274"
$ws.Range("E44").Value = "  +1.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.15"
$ws.Range("E45").Value = "  +3.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1200"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.594"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.41"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05534"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.364"
$ws.Range("E50").Value = "  -3.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3686"
$ws.Range("E51").Value = "  -0.80%  "
